# Apply updated crypto price / 1h-volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.561.72"
$ws.Range("E2").Value = "  +1.45%  "

$ws.Range("D3").Value = "3.088.44"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.51%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "3.086.38"
$ws.Range("E8").Value = "  +0.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.44%  "

$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("E12").Value = "  +4.16%  "

$ws.Range("E13").Value = "  +1.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.46%  "

$ws.Range("D15").Value = "3.588.38"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").Value = "64.569.44"
$ws.Range("E16").Value = "  +1.41%  "

$ws.Range("D17").Value = "3.089.49"
$ws.Range("E17").Value = "  +0.75%  "

$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.81%  "

$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.53%  "

$ws.Range("E24").Value = "  +8.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  +1.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.84%  "

$ws.Range("E29").Value = "  +4.17%  "

$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.29%  "

$ws.Range("E35").Value = "  +3.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "470.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.90%  "

$ws.Range("E38").Value = "  +19.78%  "

$ws.Range("E39").Value = "  +3.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0408"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.10%  "

$ws.Range("D41").Value = "2.983.32"
$ws.Range("E41").Value = "  -5.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.12%  "

$ws.Range("E43").Value = "  -3.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("E45").Value = "  +3.89%  "

$ws.Range("E46").Value = "  +6.36%  "

$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("E48").Value = "  +2.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.96%  "

$ws.Range("D50").Value = "0.0₃0524"
$ws.Range("E50").Value = "  +2.65%  "

$ws.Range("E51").Value = "  +1.06%  "
